# Minor updates on slides:
#  - Slide 8 : tweak wording "delay to the first" -> "delay until the first"
#  - Slide 9 : reword first bullet and drop the trailing "(Can also move ...)" sentence
#  - Slide 9 : nudge two diagram rectangles vertically (BR Switch diagram)
#  - Slide 10: nudge two diagram rectangles vertically (CPE Switch diagram)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 8: "(Can also delay to the first packet_in from the CPE)"
#       -> "(Can also delay until the first packet_in from the CPE)"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$sh8 = $s8.Shapes.Item(1)
$tr8 = $sh8.TextFrame.TextRange
$t8 = $tr8.Text
$search8 = "(Can also delay to the first "
$idx8 = $t8.IndexOf($search8)
if ($idx8 -ge 0) {
    $sub8 = $tr8.Characters($idx8 + 1, $search8.Length)
    $sub8.Text = "(Can also delay until the first "
}

# ---------------------------------------------------------------------------
# Slide 9: "Controller installs forwarding rules in BR Switch (per-subscriber)
#           (Can also move to the first packet_in from the CPE)"
#       -> "For every binding state: Controller installs forwarding rules in
#           BR Switch (per-subscriber)"
#          (keep the line break, drop the trailing sentence after it)
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$sh9 = $s9.Shapes.Item(1)
$tr9 = $sh9.TextFrame.TextRange
$t9 = $tr9.Text
$breakChar = [char]11
$breakIdx = $t9.IndexOf($breakChar)
if ($breakIdx -ge 0) {
    $tailStart = $breakIdx + 2
    $tailLen = $t9.Length - $tailStart + 1
    if ($tailLen -gt 0) {
        $tail = $tr9.Characters($tailStart, $tailLen)
        $tail.Delete()
    }
}
$r9 = $tr9.Runs(1)
$r9.Text = "For every binding state: " + $r9.Text

# ---------------------------------------------------------------------------
# Slide 9: shift two rectangles ("push IP-IP6 header..." / "pop IP-IP6
# header...") up by a couple EMU (~ -0.019pt), leaving x/width/height as-is.
# (Literal point values are chosen so the Single-precision round-trip that
# PowerPoint's COM layer performs on Shape.Top lands on the exact target
# EMU offset.)
# ---------------------------------------------------------------------------
$sh9_75 = $s9.Shapes.Item(12)
$sh9_75.Top = 388.60615543228346

$sh9_24 = $s9.Shapes.Item(14)
$sh9_24.Top = 298.6122047244094

# ---------------------------------------------------------------------------
# Slide 10: shift two rectangles ("pop IP-IP6 header..." / "Packet_in") up
# slightly, leaving x/width/height as-is.
# ---------------------------------------------------------------------------
$s10 = $p.Slides.Item(10)

$sh10_83 = $s10.Shapes.Item(20)
$sh10_83.Top = 467.5765354330709

$sh10_86 = $s10.Shapes.Item(22)
$sh10_86.Top = 368.4020538440945
